# "cierre 19 May 22" -- closing out the April credits workbook.
#
# The April sheet ("Hoja3") already carries all of April's data; this pass
# only finishes dressing it up the way the other closed months
# (ENERO/FEBRERO/MARZO) already look, and leaves the cursor where the user
# left it when they saved: back on the April sheet, having just glanced at
# row 59 of March.

$wb = $excel.ActiveWorkbook

# 1) Give the April sheet its real name instead of the default "Hoja3",
#    and color its tab the same green used to mark a finished month.
$wsAbril = $wb.Worksheets.Item("Hoja3")
$wsAbril.Name = "REMISIONES   ABRIL  2022"
$wsAbril.Tab.Color = 0x50B000

# 2) Revisit March briefly (selection lands on D59) ...
$wsMarzo = $wb.Worksheets.Item("REMISIONES   MARZO   2022  ")
$wsMarzo.Select()
$wsMarzo.Range("D59").Select()

# 3) ... then come back to April, which stays the active/visible tab, with
#    the selection resting on E49.
$wsAbril.Select()
$wsAbril.Range("E49").Select()
